# Update the "Generate Report for Handback" timestamps.
# These cells are formatted with a date/time number format but the stored
# value is literal text (shared string), so we must force text assignment.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 15:32:54"

# --- zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for first file row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 15:32:49"
$wsZhCn.Range("K2").Value = "2016-09-01 15:33:23"

# --- de-de sheet: "Correspond Handback DateTime" for first file row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-01 15:33:31"
